$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore default gridline/heading display (drops the legacy showGridLines/showRowColHeaders markup)
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# New date-formula cells in B1/B2
$ws.Range("B1").NumberFormat = "yyyy-mm-dd;@"
$ws.Range("B1").Formula = "=DATE(2015,2,2)"

$ws.Range("B2").NumberFormat = "yyyy-mm-dd;@"
$ws.Range("B2").Formula = "=EOMONTH(B1,1)"

# Drop the old SUM formula in row 12
$ws.Range("A12").ClearContents()

# Column B sized to fit its new date content
$ws.Columns.Item(2).EntireColumn.AutoFit()

# Selection moves to B2
$ws.Range("B2").Select() | Out-Null
